$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 554.75
$ws.Range("I19").Value = 424.75
$ws.Range("J19").Value = 684.75
$ws.Range("K19").Value = 424.75
$ws.Range("L19").Value = 684.75
$ws.Range("M19").Value = -249.75
$ws.Range("N19").Value = -1034.75
$ws.Range("H28").Value = 1777.8572
$ws.Range("I28").Value = 1439
$ws.Range("K28").Value = 1439
$ws.Range("M28").Value = -954
$ws.Range("H104").Value = 375
$ws.Range("I104").Value = 375
$ws.Range("K104").Value = 1125
$ws.Range("M104").Value = 622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12580
$ws.Range("I74").Value = 950
$ws.Range("J74").Value = 20333.334
$ws.Range("K74").Value = 950
$ws.Range("L74").Value = 20333.334
$ws.Range("M74").Value = -76
$ws.Range("N74").Value = -22081.334
$ws.Range("H77").Value = 12580
$ws.Range("I77").Value = 950
$ws.Range("J77").Value = 20333.334
$ws.Range("K77").Value = 4750
$ws.Range("L77").Value = 101666.67
$ws.Range("M77").Value = -382
$ws.Range("N77").Value = -110402.67
$ws.Range("H92").Value = 54998
$ws.Range("J92").Value = 54998
$ws.Range("L92").Value = 54998
$ws.Range("N92").Value = -59990
$ws.Range("H125").Value = 49998.5
$ws.Range("J125").Value = 49998.5
$ws.Range("L125").Value = 49998.5
$ws.Range("N125").Value = -59838.5
$ws.Range("H132").Value = 14102.444
$ws.Range("I132").Value = 12417.429
$ws.Range("K132").Value = 37252.287
$ws.Range("M132").Value = -34722.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 26997.334
$ws.Range("J76").Value = 26996.5
$ws.Range("L76").Value = 26996.5
$ws.Range("N76").Value = -27626.5
$ws.Range("H79").Value = 26997.334
$ws.Range("J79").Value = 26996.5
$ws.Range("L79").Value = 26996.5
$ws.Range("N79").Value = -29180.5
$ws.Range("H80").Value = 1982.375
$ws.Range("I80").Value = 1812.4
$ws.Range("J80").Value = 2265.6667
$ws.Range("K80").Value = 1812.4
$ws.Range("L80").Value = 2265.6667
$ws.Range("M80").Value = -814.4000000000001
$ws.Range("N80").Value = -4261.6667
$ws.Range("H83").Value = 1982.375
$ws.Range("I83").Value = 1812.4
$ws.Range("J83").Value = 2265.6667
$ws.Range("K83").Value = 9062
$ws.Range("L83").Value = 11328.3335
$ws.Range("M83").Value = -4070
$ws.Range("N83").Value = -21312.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.85715
$ws.Range("I7").Value = 86.25
$ws.Range("J7").Value = 230
$ws.Range("K7").Value = 86.25
$ws.Range("L7").Value = 230
$ws.Range("M7").Value = 26.75
$ws.Range("N7").Value = -456
$ws.Range("H16").Value = 2132.6667
$ws.Range("J16").Value = 399
$ws.Range("L16").Value = 399
$ws.Range("N16").Value = -973
$ws.Range("H92").Value = 31926.125
$ws.Range("J92").Value = 31926.125
$ws.Range("L92").Value = 31926.125
$ws.Range("N92").Value = -36918.125
$ws.Range("H99").Value = 2000
$ws.Range("J99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 2498.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2498.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2498.5
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -5992.5
$ws.Range("H107").Value = 629.0625
$ws.Range("I107").Value = 685.1429000000001
$ws.Range("J107").Value = 236.5
$ws.Range("K107").Value = 685.1429000000001
$ws.Range("L107").Value = 236.5
$ws.Range("M107").Value = 1234.8571
$ws.Range("N107").Value = -4076.5
$ws.Range("H113").Value = 2132.6667
$ws.Range("J113").Value = 399
$ws.Range("L113").Value = 399
$ws.Range("N113").Value = -4739
$ws.Range("H122").Value = 1399.3334
$ws.Range("J122").Value = 1500
$ws.Range("L122").Value = 4500
$ws.Range("N122").Value = -9400
$ws.Range("H125").Value = 49995
$ws.Range("J125").Value = 49995
$ws.Range("L125").Value = 49995
$ws.Range("N125").Value = -54915
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940
$ws.Range("H134").Value = 3779.818
$ws.Range("I134").Value = 1697.25
$ws.Range("J134").Value = 9333.333000000001
$ws.Range("K134").Value = 5091.75
$ws.Range("L134").Value = 27999.999
$ws.Range("M134").Value = -2556.75
$ws.Range("N134").Value = -33069.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16080
$ws.Range("H131").Value = 2161.8462
$ws.Range("I131").Value = 4499.25
$ws.Range("K131").Value = 13497.75
$ws.Range("M131").Value = -8457.75
$ws.Range("H134").Value = 8998.75
$ws.Range("I134").Value = 10998.5
$ws.Range("K134").Value = 32995.5
$ws.Range("M134").Value = -27925.5
$ws.Range("H139").Value = 1479.2
$ws.Range("I139").Value = 1479.2
$ws.Range("K139").Value = 4437.6
$ws.Range("M139").Value = 702.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1852.5
$ws.Range("I31").Value = 2112.5
$ws.Range("J31").Value = 1592.5
$ws.Range("K31").Value = 2112.5
$ws.Range("L31").Value = 1592.5
$ws.Range("M31").Value = -1820.5
$ws.Range("N31").Value = -2176.5
$ws.Range("H37").Value = 1852.5
$ws.Range("I37").Value = 2112.5
$ws.Range("J37").Value = 1592.5
$ws.Range("K37").Value = 2112.5
$ws.Range("L37").Value = 1592.5
$ws.Range("M37").Value = -1835.5
$ws.Range("N37").Value = -2146.5
$ws.Range("H92").Value = 7125.8335
$ws.Range("J92").Value = 7125.8335
$ws.Range("L92").Value = 7125.8335
$ws.Range("N92").Value = -10869.8335
$ws.Range("H132").Value = 9536.571
$ws.Range("I132").Value = 7939.25
$ws.Range("K132").Value = 23817.75
$ws.Range("M132").Value = -21287.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2183.1667
$ws.Range("J22").Value = 5399.5
$ws.Range("L22").Value = 5399.5
$ws.Range("N22").Value = -5989.5
$ws.Range("H27").Value = 2183.1667
$ws.Range("J27").Value = 5399.5
$ws.Range("L27").Value = 5399.5
$ws.Range("N27").Value = -5613.5
$ws.Range("H46").Value = 4654.778
$ws.Range("J46").Value = 5756.2856
$ws.Range("L46").Value = 5756.2856
$ws.Range("N46").Value = -6132.2856
$ws.Range("H110").Value = 29329.6
$ws.Range("J110").Value = 29329.6
$ws.Range("L110").Value = 29329.6
$ws.Range("N110").Value = -37509.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 26502500
$ws.Range("I4").Value = 26502500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 26502500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -26502387
$ws.Range("N4").ClearContents()
$ws.Range("H132").Value = 8203.333000000001
$ws.Range("I132").Value = 5666
$ws.Range("J132").Value = 10740.667
$ws.Range("K132").Value = 16998
$ws.Range("L132").Value = 32222.001
$ws.Range("M132").Value = -14468
$ws.Range("N132").Value = -37282.001
$ws.Range("H136").Value = 13000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450
